# Update the "Training Dashboard" sheet with the new progress as of 04-Nov-2025:
# - PERIOD TO EXPIRE (column H) decreases by 1 for each data row (one day closer to expiry)
# - LAST UPDATE (column I) changes from 03-Nov-2025 to 04-Nov-2025 (kept as literal text,
#   not converted to a real date, matching the original inline-string cell content)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow = 28

# --- Column H: PERIOD TO EXPIRE -> decrement each value by 1 ---
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)   # Column H
    $periodCell.Value2 = $periodCell.Value2 - 1
}

# --- Column I: LAST UPDATE -> "04-Nov-2025" (plain text, not a date value) ---
$lastUpdateRange = $ws.Range("I$firstRow`:I$lastRow")
$lastUpdateRange.Formula = "=""04-Nov-2025"""
$lastUpdateRange.Copy($lastUpdateRange)
$lastUpdateRange.PasteSpecial(-4163)  # xlPasteValues - flattens formula to a literal value
$excel.CutCopyMode = $false
